$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '256.47'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-0.14%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '26.76'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-0.32%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '4.652'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-1.75%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05949'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '0.49%'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-0.69%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8507'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-2.02%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9088'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-3.80%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1377'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-1.96%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.04559'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '21.39%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07000'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-1.66%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03055'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-3.54%'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-1.84%'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-1.16%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0006069'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.19%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.006019'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-1.16%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.467'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-0.96%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.157'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-1.40%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.3028'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-3.68%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1295'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '1.42%'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '1.52%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04241'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '0.64%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001215'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-0.23%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004764'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '11.07%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001200'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '76.53%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0001523'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '2.06%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03769'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-1.43%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006198'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '56.96%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1094'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-0.71%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002320'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '7.25%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01451'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '29.97%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005257'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-4.38%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000750'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '0.02%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.04299'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-51.41%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.2415'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '9,798.44%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002099'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '0.02%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0001999'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.02%'
